# Change network line settings:
# The R + L + (C // G) branch model has been changed to a (R + L) // C // G
# branch model. This only affects the "NetworkLine" sheet (user-data table)
# and the descriptive notes on both "NetworkLine" and "NetworkLine_IEEE".

$wb = $excel.ActiveWorkbook

$wsLine = $wb.Worksheets.Item("NetworkLine")
$wsIEEE = $wb.Worksheets.Item("NetworkLine_IEEE")

# Update the explanatory note describing the branch model on each sheet.
$wsLine.Range("A3").Value = "In this form, Ybranch = 1/(R+jwL)+(G+jwC), i.e., R and L are in series, G and C are in prallel, RL and GC are in parallel."
$wsIEEE.Range("A3").Value = "In this form, a pi-circuit between two buses, i.e., series R+jwL impedance with parallel G/2+jwC/2 admittance between each bus and ground."

# The "G (pu)" column for the mutual (line) branches used to hold "inf" and
# now holds a plain 0 (no shunt conductance on those rows).
$wsLine.Range("F11:F14").Value = 0

# The "R (pu)" / "wL (pu)" columns for the self branches used to hold 0 and
# now hold "inf" (i.e. the series RL leg of those rows is now open / removed).
$wsLine.Range("C15:D18").Value = "inf"

# Restore per-sheet selections and make "NetworkLine" the active tab.
$wb.Worksheets.Item("Bus").Range("D14").Select() | Out-Null
$wb.Worksheets.Item("Apparatus").Range("A1").Select() | Out-Null
$wsIEEE.Range("D8").Select() | Out-Null
$wb.Worksheets.Item("Basic").Range("A2").Select() | Out-Null
$wb.Worksheets.Item("Advance").Range("B16").Select() | Out-Null

$wsLine.Activate() | Out-Null
$wsLine.Range("D15").Select() | Out-Null
